# Helper: force a value to be stored as TEXT (not auto-coerced to a number),
# without leaving any stray NumberFormat/Style behind. We do this by writing
# a formula that evaluates to the literal text, then converting that formula
# to a static value via Copy + PasteSpecial(xlPasteValues).
function Set-TextValue {
    param($Range, $Text)
    $Range.Formula = '="' + $Text + '"'
    $Range.Copy()
    $Range.PasteSpecial(-4163)
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q4" sheet right after "总计" (i.e. before the
#    existing "2022-Q1" sheet), matching the new tab order from the diff.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# NOTE: sheet references obtained via Worksheets.Item(...) are resolved by
# *position*, not object identity - they go stale (silently point at a
# different sheet) once the tab order changes underneath them. So we only
# look up "2022-Q1" by name AFTER the insert above has finished reshuffling
# the tab order.
$q1Sheet = $wb.Worksheets.Item("2022-Q1")

# Clone the header row (and its styling) from the "2022-Q1" sheet so the new
# sheet matches the look of its neighbours, then stamp the 2022-Q4 fund data.
# (Column A of the header row is intentionally left untouched/blank, same as
# every other sheet in this workbook - only B1:H1 carry header labels.)
$q1Sheet.Range("B1:H1").Copy($newSheet.Range("B1"))
$q1Sheet.Range("A2:H2").Copy($newSheet.Range("A2"))
$q1Sheet.Range("A2:H2").Copy($newSheet.Range("A3"))
$q1Sheet.Range("A2:H2").Copy($newSheet.Range("A4"))
$q1Sheet.Range("A2:H2").Copy($newSheet.Range("A5"))

$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "005413"
Set-TextValue $newSheet.Range("C2") "金信民长灵活配置混合C"
Set-TextValue $newSheet.Range("D2") "0.86"
Set-TextValue $newSheet.Range("E2") "89.93"
Set-TextValue $newSheet.Range("F2") "4.47"
Set-TextValue $newSheet.Range("G2") "0.0384"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "005412"
Set-TextValue $newSheet.Range("C3") "金信民长灵活配置混合A"
Set-TextValue $newSheet.Range("D3") "0.83"
Set-TextValue $newSheet.Range("E3") "89.93"
Set-TextValue $newSheet.Range("F3") "4.47"
Set-TextValue $newSheet.Range("G3") "0.0371"
$newSheet.Range("H3").Value = 10

$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet.Range("B4") "020034"
Set-TextValue $newSheet.Range("C4") "国泰民安增利债券C"
Set-TextValue $newSheet.Range("D4") "1.06"
Set-TextValue $newSheet.Range("E4") "49.57"
Set-TextValue $newSheet.Range("F4") "2.80"
Set-TextValue $newSheet.Range("G4") "0.0297"
$newSheet.Range("H4").Value = 7

$newSheet.Range("A5").Value = 3
Set-TextValue $newSheet.Range("B5") "020033"
Set-TextValue $newSheet.Range("C5") "国泰民安增利债券A"
Set-TextValue $newSheet.Range("D5") "0.20"
Set-TextValue $newSheet.Range("E5") "49.57"
Set-TextValue $newSheet.Range("F5") "2.80"
Set-TextValue $newSheet.Range("G5") "0.0056"
$newSheet.Range("H5").Value = 7

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: push existing rows 2-5 down to 3-6
#    (copying values + formatting so no new styles get introduced), then
#    write the brand-new row 2 for "2022-Q4", and fix up the running index
#    in column A.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")

$ws1.Range("A5:D5").Copy($ws1.Range("A6"))
$ws1.Range("A4:D4").Copy($ws1.Range("A5"))
$ws1.Range("A3:D3").Copy($ws1.Range("A4"))
$ws1.Range("A2:D2").Copy($ws1.Range("A3"))

$ws1.Range("A2").Value = 0
Set-TextValue $ws1.Range("B2") "2022-Q4"
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 0.11

$ws1.Range("A3").Value = 1
$ws1.Range("A4").Value = 2
$ws1.Range("A5").Value = 3
$ws1.Range("A6").Value = 4

Write-Host "2022-Q4 sheet added and 总计 updated"
